$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the "Fecha" (date) and measurement columns between row 2
# and row 5, while leaving the identification columns (A, B, C, E, F, G, H, O, R)
# untouched, since they are identical on both rows already.

# Row 2 (was the "Primera" / 2022-03-03 record) becomes the
# "Segunda" / 2021-06-23 record.
$ws.Range("D2").Value = 44370
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 1080
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("P2").Value = 180
$ws.Range("Q2").Value = 6

# Row 5 (was the "Segunda" / 2021-06-23 record) becomes the
# "Primera" / 2022-03-03 record.
$ws.Range("D5").Value = 44623
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("N5").Value = "$/paquete"
$ws.Range("P5").Value = 1900
$ws.Range("Q5").Value = 1
